$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename sheets to their new (upper-cased / accented) names
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
